$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7505792379379272
$ws.Range("B1").Value = 3.100512504577637
$ws.Range("C1").Value = 2.903003692626953
$ws.Range("D1").Value = 2.440034151077271
$ws.Range("E1").Value = 2.13351035118103
